$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Corrections to existing rows (19 & 21): client phone number typo fix,
#    and row 19 status flipped from "Pending" to "Done".
# ---------------------------------------------------------------------------
$ws.Range("I19").Value2 = 9108464458
$ws.Range("N19").Value2 = "Done"
$ws.Range("I21").Value2 = 9108464458

# ---------------------------------------------------------------------------
# 2) New daily-log entries appended into the second table (columns G:N) for
#    rows 22-29. Each "date" row in column G re-uses the number format /
#    alignment already used by the existing log rows, so we copy that
#    formatting across before writing the new values.
# ---------------------------------------------------------------------------

function Set-LogRow($Row, $Date, $Client, $Phone, $MeetType, $MeetingId, $Password, $Issue, $Status) {

    # Date (column G) - numeric date format, centered (style already used by
    # the column, e.g. row 18).
    $ws.Range("G18").Copy()
    $ws.Range("G$Row").PasteSpecial(-4122)
    $ws.Range("G$Row").Value2 = $Date

    if ($Client -ne $null) {
        $ws.Range("H$Row").Value2 = $Client
    }

    if ($Phone -ne $null) {
        $ws.Range("I18").Copy()
        $ws.Range("I$Row").PasteSpecial(-4122)
        $ws.Range("I$Row").Value2 = $Phone
    }

    if ($MeetType -ne $null) {
        $ws.Range("J$Row").Value2 = $MeetType
    }

    if ($MeetingId -ne $null) {
        $ws.Range("K$Row").Value2 = $MeetingId
    }

    if ($Password -ne $null) {
        $ws.Range("L$Row").Value2 = $Password
    }

    if ($Issue -ne $null) {
        $ws.Range("M$Row").Value2 = $Issue
    }

    if ($Status -ne $null) {
        $ws.Range("N18").Copy()
        $ws.Range("N$Row").PasteSpecial(-4122)
        $ws.Range("N$Row").Value2 = $Status
    }
}

# Row 22 - Annapurneswari lab
Set-LogRow 22 44996 "Annapurneswari lab" 8197816681 $null $null $null "Sign Adding" "Pending"

# Row 23 - GKVK (trailing space preserved to match existing entry)
Set-LogRow 23 44996 "GKVK " 9483491342 "Team-Viewer" 832388342 $null "Test Creation" "Done"

# Row 24 - Advaita polytechnic
Set-LogRow 24 44996 "Advaita polytechnic" 6364132225 "Anydesk" "282 493 214" "N/A" "Name change & mail" "Done"

# Row 25 - blank placeholder day (wrap-text formatted, otherwise empty)
$ws.Range("G18").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G25").Value2 = 44997
$ws.Range("H25:N25").WrapText = $true

# Row 26 - kubra nursing home
Set-LogRow 26 44998 "kubra nursing home" 9980423641 "Team-Viewer" $null $null "Re-installation" "Done"

# Row 27 - Ashwini lab
Set-LogRow 27 44999 "Ashwini lab" 9663855462 "Anydesk" 793718041 "N/A" "In Report amount not showing" "Done"

# Row 28 - Skanda Lab
Set-LogRow 28 44999 "Skanda Lab" 9886410873 "Anydesk" 484406818 "N/A" "Mail rpt header" "Done"

# Row 29 - SHS lab
Set-LogRow 29 44999 "SHS lab" 9901909968 "Anydesk" 434591205 "N/A" "Report erorr" "Done"

# ---------------------------------------------------------------------------
# 3) View state: scroll the visible window down/right and move the active
#    selection, matching where the author was working when they saved.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H32").Select()
